$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.134897470474243
$ws.Range("B1").Value = 5.91640043258667
$ws.Range("C1").Value = 2.504183292388916
$ws.Range("D1").Value = 1.143244385719299
$ws.Range("E1").Value = 0.8142873048782349
